# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) / DialogAct (col J) pairs
# for the rows whose automatic annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 6;   Tag = "ba"; Act = "Appreciation" },
    @{ Row = 7;   Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 19;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 24;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 27;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 38;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 39;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 45;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 46;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 48;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 61;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 62;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 69;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 76;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 79;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 81;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 83;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 86;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 87;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 93;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 94;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 95;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 97;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 98;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 99;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 113; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
